$d = $word.ActiveDocument
$wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'

# --- 1) Title paragraph: add spacing after=0 -----------------------------
$d.Paragraphs.Item(1).SpaceAfter = 0

# --- 2) Split "Antecedents personnels : [] <br> Antecedents familiaux : []" --
#        into two separate paragraphs, both spacing after=0
$rng = $d.Paragraphs.Item(3).Range
[void]$rng.Find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, "^p", 2)
$d.Paragraphs.Item(3).SpaceAfter = 0
$d.Paragraphs.Item(4).SpaceAfter = 0

# --- 3) Insert new empty bold paragraph after "Antecedents familiaux" -------
$d.Paragraphs.Item(4).Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item(5)
$boldEmptyXml = "<w:p $wNs><w:pPr><w:spacing w:after=`"0`"/><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr></w:p>"
[void]$newPara.Range.InsertXML($boldEmptyXml)

# --- 4) Split the "Mammographie..." paragraph (index 7) into three ----------
#        paragraphs at the two manual line breaks, spacing after=0 on each,
#        and retext the first one to "[Incidences]"
$rng7 = $d.Paragraphs.Item(7).Range
[void]$rng7.Find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, "^p", 2)
$rng7b = $d.Paragraphs.Item(7).Range
[void]$rng7b.Find.Execute("^l", $false, $false, $false, $false, $false, $true, 1, $false, "^p", 2)
$d.Paragraphs.Item(7).SpaceAfter = 0
$d.Paragraphs.Item(8).SpaceAfter = 0
$d.Paragraphs.Item(9).SpaceAfter = 0

# retext "Mammographie bilatérale (face + oblique)" -> "[Incidences]"
$rngInc = $d.Paragraphs.Item(7).Range
[void]$rngInc.Find.Execute("Mammographie bilatérale (face + oblique)", $true, $false, $false, $false, $false, $true, 1, $false, "[Incidences]", 2)

# --- 5) "Échographie ..." paragraph (index 10): add spacing after=0 --------
$d.Paragraphs.Item(10).SpaceAfter = 0

# --- 6) Insert new empty bold paragraph after "Échographie ..." ------------
$d.Paragraphs.Item(10).Range.InsertParagraphAfter()
$newPara2 = $d.Paragraphs.Item(11)
[void]$newPara2.Range.InsertXML($boldEmptyXml)

# --- 7) "CONCLUSION : ..." paragraph (index 12): add spacing after=0 -------
$d.Paragraphs.Item(12).SpaceAfter = 0

# --- 8) Remove the trailing empty paragraph at the end of the document -----
$n = $d.Paragraphs.Count
$pPrev = $d.Paragraphs.Item($n - 1)
$pLast = $d.Paragraphs.Item($n)
$delRange = $d.Range($pPrev.Range.End - 1, $pLast.Range.End)
$delRange.Delete()

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    Write-Host "$($i): [$($d.Paragraphs.Item($i).Range.Text)]"
}
